# Update "Data" table: job ads (Job adverts by occupation) data refresh.
# The raw data source format changed (no longer includes 2-digit SOC, LEP/LSIP/MCA
# summaries), so the "latest period" / "next period" release labels for the
# "Job adverts by occupation" row move on by one month, and the "Highest
# qualification level by age and gender" row's "next period" label is
# re-pointed at the (unchanged) existing label text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 = "Highest qualification level by age and gender": next period text unchanged,
# just re-set to keep it pointing at the right shared string after the table shuffle.
$ws.Range("D7").Value = "Jan 2025 - Dec 2025 (Apr 26)"

# Row 13 = "Job adverts by occupation": latest/next period values roll forward one month.
$ws.Range("C13").Value = "May 2025 (26/06/25)"
$ws.Range("D13").Value = "Jun 2025 (Jul 2025)"

# Move the active selection to D14, matching the saved workbook view state.
$ws.Range("D14").Select()
